# "Generate Report for Handoff"
#
# For each localized-language sheet (zh-cn, de-de), the rows that were
# "Ready for handoff" (Priority = "low", i.e. not yet handed off) just had
# a handoff xliff generated for them: their Priority flips from "low" to
# "ht" and their Latest Handoff Datetime is refreshed to the generation
# timestamp for that language.

$wb = $excel.ActiveWorkbook

$rows = 4, 5, 6, 7

$sheetHandoffTimes = @{
    "zh-cn" = "2016-08-30 00:33:26"
    "de-de" = "2016-08-30 00:33:31"
}

foreach ($sheetName in "zh-cn", "de-de") {
    $ws = $wb.Worksheets.Item($sheetName)
    $handoffTime = $sheetHandoffTimes[$sheetName]
    foreach ($r in $rows) {
        $ws.Range("E$r").Value = "ht"
        $ws.Range("H$r").Value = $handoffTime
    }
}

# The Overview sheet's "Latest HO Xliff Generate Date" column shared the
# same underlying string as the de-de sheet's handoff datetime for these
# rows, so it picks up the refreshed de-de timestamp as well.
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Range("G$r").Value = $sheetHandoffTimes["de-de"]
}
